$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 191; this shifts the existing rows
# 191-220 down to 192-221 (matching the rest of the diff, where every
# row N (192<=N<=221) now holds what used to be row N-1's data).
$ws.Rows.Item(191).Insert()

# Populate the newly inserted row 191 with the new weekly price entry.
$ws.Range("A191").Value = 10
$ws.Range("B191").Value = 'Vega Modelo de Temuco'
$ws.Range("C191").Value = 'La Araucanía'
$ws.Range("D191").Value = 44505
$ws.Range("E191").Value = 9
$ws.Range("F191").Value = 100114013
$ws.Range("G191").Value = 'Zanahoria'
$ws.Range("H191").Value = 'Sin especificar'
$ws.Range("I191").Value = 'Primera'
$ws.Range("J191").Value = 65
$ws.Range("K191").Value = 7000
$ws.Range("L191").Value = 7000
$ws.Range("M191").Value = 7000
$ws.Range("N191").Value = '$/saco 20 kilos'
$ws.Range("O191").Value = 'Región del Maule'
$ws.Range("P191").Value = 350
$ws.Range("Q191").Value = 20
$ws.Range("R191").Value = 'Hortaliza'
